$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Helper data for the two new localized files being handed off:
#   8445301b-c14c-420c-b989-0631f4d6c5a9  (handoff 2016-03-03 07:18:18 / 07:18:28)
#   b01c9cac-feae-4ec9-aa79-354b047dc00c  (same handoff batch)
# ---------------------------------------------------------------

$mdUrlBase   = "https://github.com/OpenLocalizationTest/oltest/blob/25b6b799517af5342f052b2a1faf278ca6b4e250/e2e/"
$cfgUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/25b6b799517af5342f052b2a1faf278ca6b4e250/.localization-config"
$zhUrlBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2aed74ceb867b09ff7874b1f9ade85dd28ef7bfb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$deUrlBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0eaed73367c6511d37ce46c4e5112e4ed3e823ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$file1 = "8445301b-c14c-420c-b989-0631f4d6c5a9"
$file2 = "b01c9cac-feae-4ec9-aa79-354b047dc00c"

$file1Md  = $file1 + ".md"
$file2Md  = $file2 + ".md"
$file1Zh  = $file1 + ".ad60c060228846307d07b5707ed60c63e56c6f51.zh-cn.xlf"
$file2Zh  = $file2 + ".c545d91c876b29ff12095a9dc2ad69feceec2931.zh-cn.xlf"
$file1De  = $file1 + ".ad60c060228846307d07b5707ed60c63e56c6f51.de-de.xlf"
$file2De  = $file2 + ".c545d91c876b29ff12095a9dc2ad69feceec2931.de-de.xlf"

$handoffZh = "2016-03-03 07:18:18"
$handoffDe = "2016-03-03 07:18:28"
$epoch     = "0001-01-01 00:00:00"

# =================================================================
# Sheet "Overview": File Name | zh-cn | de-de
# =================================================================
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = $file1Md
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = $file2Md
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + "5d796f67-eac0-4609-b008-58aeb6e7f994.md"), "", "", "5d796f67-eac0-4609-b008-58aeb6e7f994.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + "92623a59-9c3e-48fc-8a69-c069db720aa6.md"), "", "", "92623a59-9c3e-48fc-8a69-c069db720aa6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($mdUrlBase + $file1Md), "", "", $file1Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($mdUrlBase + $file2Md), "", "", $file2Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config") | Out-Null

# =================================================================
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# =================================================================
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = $file1Md
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $file1Zh
$ws.Range("D4").Value = $handoffZh
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Include"

$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = $file2Md
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = $file2Zh
$ws.Range("D5").Value = $handoffZh
$ws.Range("G5").Value = $epoch
$ws.Range("H5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + "5d796f67-eac0-4609-b008-58aeb6e7f994.md"), "", "", "5d796f67-eac0-4609-b008-58aeb6e7f994.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($zhUrlBase + "5d796f67-eac0-4609-b008-58aeb6e7f994.e0c0307125cca1aaf68ef6b22b5ae4c13f22838c.zh-cn.xlf"), "", "", "5d796f67-eac0-4609-b008-58aeb6e7f994.e0c0307125cca1aaf68ef6b22b5ae4c13f22838c.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + "92623a59-9c3e-48fc-8a69-c069db720aa6.md"), "", "", "92623a59-9c3e-48fc-8a69-c069db720aa6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($zhUrlBase + "92623a59-9c3e-48fc-8a69-c069db720aa6.7c771d0bd993e919a6b46d8e46b8a34ebc6f5569.zh-cn.xlf"), "", "", "92623a59-9c3e-48fc-8a69-c069db720aa6.7c771d0bd993e919a6b46d8e46b8a34ebc6f5569.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($mdUrlBase + $file1Md), "", "", $file1Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), ($zhUrlBase + $file1Zh), "", "", $file1Zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($mdUrlBase + $file2Md), "", "", $file2Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), ($zhUrlBase + $file2Zh), "", "", $file2Zh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config") | Out-Null

# =================================================================
# Sheet "de-de": same columns, de-de xlf/datetime variants
# =================================================================
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = $file1Md
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = $file1De
$ws.Range("D4").Value = $handoffDe
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Include"

$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = $file2Md
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = $file2De
$ws.Range("D5").Value = $handoffDe
$ws.Range("G5").Value = $epoch
$ws.Range("H5").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), ($mdUrlBase + "5d796f67-eac0-4609-b008-58aeb6e7f994.md"), "", "", "5d796f67-eac0-4609-b008-58aeb6e7f994.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), ($deUrlBase + "5d796f67-eac0-4609-b008-58aeb6e7f994.e0c0307125cca1aaf68ef6b22b5ae4c13f22838c.de-de.xlf"), "", "", "5d796f67-eac0-4609-b008-58aeb6e7f994.e0c0307125cca1aaf68ef6b22b5ae4c13f22838c.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), ($mdUrlBase + "92623a59-9c3e-48fc-8a69-c069db720aa6.md"), "", "", "92623a59-9c3e-48fc-8a69-c069db720aa6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), ($deUrlBase + "92623a59-9c3e-48fc-8a69-c069db720aa6.7c771d0bd993e919a6b46d8e46b8a34ebc6f5569.de-de.xlf"), "", "", "92623a59-9c3e-48fc-8a69-c069db720aa6.7c771d0bd993e919a6b46d8e46b8a34ebc6f5569.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), ($mdUrlBase + $file1Md), "", "", $file1Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), ($deUrlBase + $file1De), "", "", $file1De) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), ($mdUrlBase + $file2Md), "", "", $file2Md) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), ($deUrlBase + $file2De), "", "", $file2De) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), $cfgUrl, "", "", ".localization-config") | Out-Null

Write-Output "Report generated for handoff"
